# Update specific imputed values in result_data_RandomForest.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"  = -12.00039999999999
    "B9"  = 6.400599999999995
    "C9"  = -11.0217
    "B18" = 6.251899999999997
    "B20" = 9.522799999999988
    "C23" = -12.1414
    "C24" = -13.3882
    "C26" = -12.7496
    "B27" = 6.456400000000006
    "C34" = -11.89680000000001
    "C35" = -12.32910000000001
    "C48" = -11.7058
    "C52" = -11.0427
    "C66" = -11.05560000000001
    "C67" = -11.3096
    "B69" = 5.138899999999994
    "B76" = 5.534999999999999
    "C80" = -13.12380000000001
    "B82" = 5.477000000000001
    "C99" = -12.6589
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
